$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '42.734.85'
$ws.Range('E2').Value = '  +1.34%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.560.97'
$ws.Range('E3').Value = '  +2.30%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.997'
$ws.Range('E4').Value = '  -0.42%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '313.44'
$ws.Range('E5').Value = '  -0.04%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '101.09'
$ws.Range('E6').Value = '  +7.18%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.570'
$ws.Range('E7').Value = '  +0.42%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.997'
$ws.Range('E8').Value = '  -0.45%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.531'
$ws.Range('E9').Value = '  +0.85%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '36.41'
$ws.Range('E10').Value = '  +4.45%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0805'
$ws.Range('E11').Value = '  +0.68%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '7.43'
$ws.Range('E12').Value = '  +1.06%  '

$ws.Range('E13').Value = '  +0.73%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '2.934.03'
$ws.Range('E14').Value = '  +1.27%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '16.13'
$ws.Range('E15').Value = '  +9.07%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '2.588.16'
$ws.Range('E16').Value = '  +2.42%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.848'
$ws.Range('E17').Value = '  +2.22%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '42.690.99'
$ws.Range('E18').Value = '  +0.56%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.81'
$ws.Range('E19').Value = '  +1.50%  '

$ws.Range('B20').Value = 'InternetComputer(DFINITY)'
$ws.Range('C20').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '12.41'
$ws.Range('E20').Value = '  +1.00%  '

$ws.Range('B21').Value = 'ShibaInu'
$ws.Range('C21').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.0₃0956'
$ws.Range('E21').Value = '  +1.16%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '69.23'
$ws.Range('E22').Value = '  +0.69%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '245.27'
$ws.Range('E23').Value = '  -1.81%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.91'
$ws.Range('E24').Value = '  +0.06%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.07'
$ws.Range('E25').Value = '  +2.86%  '

$ws.Range('B26').Value = 'EthereumClassic'
$ws.Range('C26').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '26.51'
$ws.Range('E26').Value = '  +1.15%  '

$ws.Range('B27').Value = 'Dai'
$ws.Range('C27').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.00'
$ws.Range('E27').Value = '  -0.15%  '

$ws.Range('B28').Value = 'InjectiveProtocol'
$ws.Range('C28').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '40.68'
$ws.Range('E28').Value = '  +2.83%  '

$ws.Range('B29').Value = 'Toncoin'
$ws.Range('C29').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.36'
$ws.Range('E29').Value = '  -0.93%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '10.15'
$ws.Range('E30').Value = '  +0.75%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '157.74'
$ws.Range('E31').Value = '  +1.68%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '5.72'
$ws.Range('E32').Value = '  -0.45%  '

$ws.Range('E33').Value = '  +13.99%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.0805'
$ws.Range('E34').Value = '  +3.14%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.08'
$ws.Range('E35').Value = '  +1.37%  '

$ws.Range('E36').Value = '  -2.71%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '3.21'
$ws.Range('E37').Value = '  -1.06%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '18.32'
$ws.Range('E38').Value = '  -3.35%  '

$ws.Range('E39').Value = '  +0.44%  '

$ws.Range('E40').Value = '  +1.34%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '4.22'
$ws.Range('E41').Value = '  +13.45%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '21.83'
$ws.Range('E42').Value = '  -0.06%  '

$ws.Range('B43').Value = 'NEARProtocol'
$ws.Range('C43').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '3.33'
$ws.Range('E43').Value = '  +4.69%  '

$ws.Range('B44').Value = 'FirstDigitalUSD'
$ws.Range('C44').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.00'
$ws.Range('E44').Value = '  -0.18%  '

$ws.Range('E45').Value = '  -0.16%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.969.64'
$ws.Range('E46').Value = '  -0.24%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '8.91'
$ws.Range('E47').Value = '  -0.15%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.791.22'
$ws.Range('E48').Value = '  +0.87%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '81.29'
$ws.Range('E49').Value = '  -2.08%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.194'
$ws.Range('E50').Value = '  +3.21%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '73.59'
$ws.Range('E51').Value = '  +1.19%  '
